# AssayTestFile.xlsx edit — "add study xlsx file io"
#
# The source diff renames the third worksheet ("Investigation" -> "Assay")
# and switches the selected/active tab from the first worksheet
# ("GreatAssay") to that renamed third worksheet (workbook.xml's
# bookViews/workbookView activeTab goes to 2, and tabSelected="1" moves
# from sheet1's sheetView to sheet3's sheetView). Reproduce both through
# the Excel object model: rename the sheet, then activate it so
# tabSelected/activeTab follow automatically, matching stock Excel
# behavior when a user clicks a different tab.

$wb = $excel.ActiveWorkbook

$assay = $wb.Worksheets.Item("Investigation")
$assay.Name = "Assay"

$assay.Activate()
